$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table previously held rows for years 2000, 2005-2013.
# Remove the rows for 2000 and 2005-2009 (rows 2-7), which shifts the
# remaining rows for 2010-2013 up to become rows 2-5.
$ws.Range("A2:E7").EntireRow.Delete()

$wb.Save()
